# SCD0012-005 - Admin SLN Melakukan edit KPI
# Rename the sheet to match the new SCD0012 numbering scheme and update the
# TC_ID cell (B2) from the placeholder Jira ref "DGS-230" to the proper test
# case id "SCD0012-005", mirroring the manual edit an admin would perform in
# the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename worksheet tab: SCD0215 -> SCD0012
$ws.Name = "SCD0012"

# Update the TC_ID value in B2
$ws.Range("B2").Value = "SCD0012-005"

# Column B was a bestFit column; with the longer text it widens automatically
$ws.Columns.Item(2).AutoFit()

# Leave the cursor on the row below the edited cell, as a user would after
# typing the value and pressing Enter
[void]$ws.Range("B3").Select()
